$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 6: student "פתחאללה חאג'312319569" carried over from last year.
# Copy the formatting of A2 (name-column style) onto A6, then set its value.
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").Value = "פתחאללה חאג'312319569"

$ws.Range("D6").Value = 40.5
$ws.Range("G6").Value = 71

# Update the selection to match the target view state
$ws.Range("F8").Select()
